$d = $word.ActiveDocument

# Locate the "Research & Data Analytics Leadership" paragraph under the
# Siege Analytics (PARTNER) role so the three new bullet points are
# inserted right after it and before the existing "Conceived, architected..."
# bullet.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Research & Data Analytics Leadership") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph 'Research & Data Analytics Leadership'"
}

$newBullets = @(
    "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters",
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

$insertAfter = $anchorIndex
foreach ($bulletText in $newBullets) {
    $anchorPara = $d.Paragraphs.Item($insertAfter)
    $anchorPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($insertAfter + 1)
    $newPara.Range.Text = $bulletText
    $insertAfter = $insertAfter + 1
}

Write-Output "Inserted $($newBullets.Count) paragraphs after index $anchorIndex"
